$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "cyano_codA_Km"
$ws.Range("B14").Value = "Wrap-Up-doc.docx"
$ws.Range("D14").Value = "cyano new notes"
$ws.Range("C14").Value = "cyano new description"

$ws.Range("C15").Select()
